# Updates cryptos list values (price + 1h volume change) per the
# Fri Dec 15 14:54:41 UTC 2023 GitHub Actions refresh, including the
# MultiversX / FraxShare row swap (rows 44-45).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.096.23'
$ws.Range('E2').Value = '  -1.18%  '
$ws.Range('D3').Value = '2.239.08'
$ws.Range('E3').Value = '  -1.88%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = "'" + '246.75'
$ws.Range('E5').Value = '  -1.47%  '
$ws.Range('D6').Value = "'" + '0.620'
$ws.Range('E6').Value = '  -2.03%  '
$ws.Range('D7').Value = "'" + '75.72'
$ws.Range('E7').Value = '  +5.28%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  -1.80%  '
$ws.Range('D10').Value = "'" + '40.33'
$ws.Range('E10').Value = '  +4.15%  '
$ws.Range('D11').Value = "'" + '0.0947'
$ws.Range('E11').Value = '  -1.95%  '
$ws.Range('D12').Value = "'" + '7.22'
$ws.Range('E12').Value = '  -1.39%  '
$ws.Range('E13').Value = '  -0.75%  '
$ws.Range('D14').Value = '2.575.95'
$ws.Range('E14').Value = '  -1.82%  '
$ws.Range('D15').Value = "'" + '14.87'
$ws.Range('E15').Value = '  -1.25%  '
$ws.Range('D16').Value = "'" + '0.860'
$ws.Range('E16').Value = '  -2.26%  '
$ws.Range('D17').Value = '2.236.73'
$ws.Range('E17').Value = '  -2.22%  '
$ws.Range('D18').Value = '42.132.99'
$ws.Range('E18').Value = '  -1.01%  '
$ws.Range('D19').Value = '0.0₃0976'
$ws.Range('E19').Value = '  -1.84%  '
$ws.Range('E20').Value = '  -1.96%  '
$ws.Range('D21').Value = "'" + '71.43'
$ws.Range('E21').Value = '  -1.36%  '
$ws.Range('D22').Value = "'" + '231.30'
$ws.Range('E22').Value = '  -1.61%  '
$ws.Range('D23').Value = "'" + '2.20'
$ws.Range('E23').Value = '  -3.81%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('E25').Value = '  -5.24%  '
$ws.Range('D26').Value = "'" + '11.15'
$ws.Range('E26').Value = '  -4.23%  '
$ws.Range('D27').Value = "'" + '2.32'
$ws.Range('E27').Value = '  -4.61%  '
$ws.Range('D28').Value = "'" + '7.11'
$ws.Range('E28').Value = '  +10.86%  '
$ws.Range('E29').Value = '  -1.37%  '
$ws.Range('D30').Value = "'" + '168.62'
$ws.Range('E30').Value = '  +0.62%  '
$ws.Range('D31').Value = "'" + '20.49'
$ws.Range('E31').Value = '  -2.78%  '
$ws.Range('E32').Value = '  +6.73%  '
$ws.Range('D33').Value = "'" + '32.61'
$ws.Range('E33').Value = '  +3.29%  '
$ws.Range('E34').Value = '  -6.61%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').Value = "'" + '4.53'
$ws.Range('E36').Value = '  -4.52%  '
$ws.Range('D37').Value = "'" + '4.84'
$ws.Range('E37').Value = '  +2.36%  '
$ws.Range('E38').Value = '  -3.12%  '
$ws.Range('D39').Value = "'" + '13.28'
$ws.Range('E39').Value = '  -5.22%  '
$ws.Range('D40').Value = "'" + '5.92'
$ws.Range('E40').Value = '  -1.07%  '
$ws.Range('E41').Value = '  -5.45%  '
$ws.Range('D42').Value = "'" + '116.79'
$ws.Range('E42').Value = '  +20.92%  '
$ws.Range('E43').Value = '  -4.74%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = "'" + '8.74'
$ws.Range('E44').Value = '  -5.26%  '
$ws.Range('B45').Value = 'MultiversX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D45').Value = "'" + '58.81'
$ws.Range('E45').Value = '  -4.99%  '
$ws.Range('D46').Value = "'" + '0.100'
$ws.Range('E46').Value = '  -2.93%  '
$ws.Range('D47').Value = "'" + '0.997'
$ws.Range('E47').Value = '  -0.43%  '
$ws.Range('E48').Value = '  -4.04%  '
$ws.Range('E49').Value = '  -1.68%  '
$ws.Range('D50').Value = "'" + '4.19'
$ws.Range('E50').Value = '  -14.18%  '
$ws.Range('D51').Value = "'" + '2.27'
$ws.Range('E51').Value = '  -0.42%  '
